$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.301687528568038
$ws.Cells.Item(2, 3).Value2 = 0.300747461727866
$ws.Cells.Item(2, 4).Value2 = 0.1233630491873043
$ws.Cells.Item(2, 5).Value2 = 0.4157543743867365
$ws.Cells.Item(2, 6).Value2 = 2.798705207953788
$ws.Cells.Item(2, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(3, 2).Value2 = 1.171273146918452
$ws.Cells.Item(3, 3).Value2 = 0.2623057578358043
$ws.Cells.Item(3, 4).Value2 = 0.1128718144349108
$ws.Cells.Item(3, 5).Value2 = 0.362076079561902
$ws.Cells.Item(3, 6).Value2 = 2.580460778543483
$ws.Cells.Item(3, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(4, 2).Value2 = 1.092247776592785
$ws.Cells.Item(4, 3).Value2 = 0.2388677213562289
$ws.Cells.Item(4, 4).Value2 = 0.1064563759749859
$ws.Cells.Item(4, 5).Value2 = 0.3292994308155812
$ws.Cells.Item(4, 6).Value2 = 2.447576670238021
$ws.Cells.Item(4, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(5, 2).Value2 = 1.060301755946682
$ws.Cells.Item(5, 3).Value2 = 0.229355877005446
$ws.Cells.Item(5, 4).Value2 = 0.1038477432212375
$ws.Cells.Item(5, 5).Value2 = 0.3159841952222422
$ws.Cells.Item(5, 6).Value2 = 2.393692262385002
$ws.Cells.Item(5, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(6, 2).Value2 = 1.055012507639333
$ws.Cells.Item(6, 3).Value2 = 0.2277787548360664
$ws.Cells.Item(6, 4).Value2 = 0.1034149009966683
$ws.Cells.Item(6, 5).Value2 = 0.3137755948689431
$ws.Cells.Item(6, 6).Value2 = 2.384760503242802
$ws.Cells.Item(6, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(7, 2).Value2 = 1.091815907367675
$ws.Cells.Item(7, 3).Value2 = 0.2387392848912668
$ws.Cells.Item(7, 4).Value2 = 0.1064211730808466
$ws.Cells.Item(7, 5).Value2 = 0.3291196943742705
$ws.Cells.Item(7, 6).Value2 = 2.446848904307927
$ws.Cells.Item(7, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(8, 2).Value2 = 1.256498824658706
$ws.Cells.Item(8, 3).Value2 = 0.2874570817163544
$ws.Cells.Item(8, 4).Value2 = 0.1197396920338889
$ws.Cells.Item(8, 5).Value2 = 0.3972057125687911
$ws.Cells.Item(8, 6).Value2 = 2.723213523426864
$ws.Cells.Item(8, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(9, 2).Value2 = 1.588085390464755
$ws.Cells.Item(9, 3).Value2 = 0.3844115841467897
$ws.Cells.Item(9, 4).Value2 = 0.1461065271936377
$ws.Cells.Item(9, 5).Value2 = 0.5323678854346667
$ws.Cells.Item(9, 6).Value2 = 3.274723129474808
$ws.Cells.Item(9, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(10, 2).Value2 = 1.837444416842061
$ws.Cells.Item(10, 3).Value2 = 0.456667782732552
$ws.Cells.Item(10, 4).Value2 = 0.1656892636845839
$ws.Cells.Item(10, 5).Value2 = 0.6329670626702324
$ws.Cells.Item(10, 6).Value2 = 3.686736685555786
$ws.Cells.Item(10, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(11, 2).Value2 = 1.952238773908789
$ws.Cells.Item(11, 3).Value2 = 0.4897972460817073
$ws.Cells.Item(11, 4).Value2 = 0.1746568592192546
$ws.Cells.Item(11, 5).Value2 = 0.6790791895586068
$ws.Cells.Item(11, 6).Value2 = 3.875874300919293
$ws.Cells.Item(11, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(12, 2).Value2 = 1.995912514945587
$ws.Cells.Item(12, 3).Value2 = 0.5023827667379237
$ws.Cells.Item(12, 4).Value2 = 0.1780622397111813
$ws.Cells.Item(12, 5).Value2 = 0.6965962014160993
$ws.Cells.Item(12, 6).Value2 = 3.947759482263052
$ws.Cells.Item(12, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(13, 2).Value2 = 1.986497406172077
$ws.Cells.Item(13, 3).Value2 = 0.4996704208222127
$ws.Cells.Item(13, 4).Value2 = 0.1773283890871369
$ws.Cells.Item(13, 5).Value2 = 0.6928210541511817
$ws.Cells.Item(13, 6).Value2 = 3.932265758005713
$ws.Cells.Item(13, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(14, 2).Value2 = 1.955827710953031
$ws.Cells.Item(14, 3).Value2 = 0.4908318444166184
$ws.Cells.Item(14, 4).Value2 = 0.1749368246897234
$ws.Cells.Item(14, 5).Value2 = 0.6805191853053856
$ws.Cells.Item(14, 6).Value2 = 3.881782971614427
$ws.Cells.Item(14, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(15, 2).Value2 = 1.937068397614212
$ws.Cells.Item(15, 3).Value2 = 0.4854232683167652
$ws.Cells.Item(15, 4).Value2 = 0.1734731962532408
$ws.Cells.Item(15, 5).Value2 = 0.6729912991885101
$ws.Cells.Item(15, 6).Value2 = 3.850895548665449
$ws.Cells.Item(15, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(16, 2).Value2 = 1.829970381361761
$ws.Cells.Item(16, 3).Value2 = 0.454508155220708
$ws.Cells.Item(16, 4).Value2 = 0.1651044902739898
$ws.Cells.Item(16, 5).Value2 = 0.6299609881346697
$ws.Cells.Item(16, 6).Value2 = 3.674411954279094
$ws.Cells.Item(16, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(17, 2).Value2 = 1.764623137037574
$ws.Cells.Item(17, 3).Value2 = 0.4356112001661359
$ws.Cells.Item(17, 4).Value2 = 0.1599864691876292
$ws.Cells.Item(17, 5).Value2 = 0.6036563113331539
$ws.Cells.Item(17, 6).Value2 = 3.566595007603013
$ws.Cells.Item(17, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(18, 2).Value2 = 1.727164610783916
$ws.Cells.Item(18, 3).Value2 = 0.4247664505884927
$ws.Cells.Item(18, 4).Value2 = 0.1570482393692032
$ws.Cells.Item(18, 5).Value2 = 0.5885590947900852
$ws.Cells.Item(18, 6).Value2 = 3.504741527073719
$ws.Cells.Item(18, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(19, 2).Value2 = 1.714503450146765
$ws.Cells.Item(19, 3).Value2 = 0.4210986920562618
$ws.Cells.Item(19, 4).Value2 = 0.1560543227081723
$ws.Cells.Item(19, 5).Value2 = 0.5834528703535824
$ws.Cells.Item(19, 6).Value2 = 3.483825968738444
$ws.Cells.Item(19, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(20, 2).Value2 = 1.771566211276365
$ws.Cells.Item(20, 3).Value2 = 0.4376202807796403
$ws.Cells.Item(20, 4).Value2 = 0.1605307126019682
$ws.Cells.Item(20, 5).Value2 = 0.6064530846480238
$ws.Cells.Item(20, 6).Value2 = 3.578055627605522
$ws.Cells.Item(20, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(21, 2).Value2 = 1.964830552200908
$ws.Cells.Item(21, 3).Value2 = 0.4934268375054671
$ws.Cells.Item(21, 4).Value2 = 0.175639017821851
$ws.Cells.Item(21, 5).Value2 = 0.6841309983901596
$ws.Cells.Item(21, 6).Value2 = 3.896603718956044
$ws.Cells.Item(21, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(22, 2).Value2 = 2.092330707953749
$ws.Cells.Item(22, 3).Value2 = 0.5301347953886193
$ws.Cells.Item(22, 4).Value2 = 0.1855693007810544
$ws.Cells.Item(22, 5).Value2 = 0.7352231267241507
$ws.Cells.Item(22, 6).Value2 = 4.106333638094725
$ws.Cells.Item(22, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(23, 2).Value2 = 2.024169780583463
$ws.Cells.Item(23, 3).Value2 = 0.5105206166376774
$ws.Cells.Item(23, 4).Value2 = 0.1802638442694331
$ws.Cells.Item(23, 5).Value2 = 0.7079228058401696
$ws.Cells.Item(23, 6).Value2 = 3.99425021284344
$ws.Cells.Item(23, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(24, 2).Value2 = 1.768426903915952
$ws.Cells.Item(24, 3).Value2 = 0.436711915334115
$ws.Cells.Item(24, 4).Value2 = 0.1602846471847954
$ws.Cells.Item(24, 5).Value2 = 0.6051885837940887
$ws.Cells.Item(24, 6).Value2 = 3.572873871679349
$ws.Cells.Item(24, 8).Value2 = 0.07973214163530429

$ws.Cells.Item(25, 2).Value2 = 1.497409031274401
$ws.Cells.Item(25, 3).Value2 = 0.3580147870183055
$ws.Cells.Item(25, 4).Value2 = 0.1389408687351619
$ws.Cells.Item(25, 5).Value2 = 0.4955969926831187
$ws.Cells.Item(25, 6).Value2 = 3.124400418051067
$ws.Cells.Item(25, 8).Value2 = 0.07973214163530429
